# feat: add 2022-Q4 data
#
# 1. Insert a new "2022-Q4" sheet (copied from "2022-Q3" so it inherits the
#    exact same layout/styles) right after "总计", then overwrite its values
#    with the new quarter's fund data.
# 2. Insert a new row into the "总计" summary sheet for the 2022-Q4 totals,
#    pushing the existing quarters down by one row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: new "2022-Q4" worksheet
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($null, $wb.Worksheets.Item("总计"))
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# Row 2: 000593
$q4.Range("B2").Value = "'000593"
$q4.Range("C2").Value = "易方达标普全球高端消费品指数增强（QDII）美元现汇"
$q4.Range("D2").Value = "'2.30"
$q4.Range("E2").Value = "'93.71"
$q4.Range("F2").Value = "'9.03"
$q4.Range("G2").Value = "'0.2077"
$q4.Range("H2").Value = 2

# Row 3: 005676
$q4.Range("B3").Value = "'005676"
$q4.Range("C3").Value = "易方达标普全球高端消费品指数增强C（QDII）人民币"
$q4.Range("D3").Value = "'2.30"
$q4.Range("E3").Value = "'93.71"
$q4.Range("F3").Value = "'9.03"
$q4.Range("G3").Value = "'0.2077"
$q4.Range("H3").Value = 2

# Row 4: 118002
$q4.Range("B4").Value = "'118002"
$q4.Range("C4").Value = "易方达标普全球高端消费品指数增强A（QDII）人民币"
$q4.Range("D4").Value = "'2.30"
$q4.Range("E4").Value = "'93.71"
$q4.Range("F4").Value = "'9.03"
$q4.Range("G4").Value = "'0.2077"
$q4.Range("H4").Value = 2

# ---------------------------------------------------------------------
# Step 2: "总计" sheet — insert a new row for 2022-Q4
# ---------------------------------------------------------------------
$tot = $wb.Worksheets.Item("总计")
$tot.Rows.Item(2).Insert()

# carry the formatting from the (now shifted) row below so the new row
# matches the rest of the table exactly
$tot.Range("A3:D3").Copy()
$tot.Range("A2:D2").PasteSpecial(-4122)

$tot.Range("A2").Value = 0
$tot.Range("B2").Value = "2022-Q4"
$tot.Range("C2").Value = 3
$tot.Range("D2").Value = 0.62

# re-index column A (0-based row index) and refresh the period labels
$tot.Range("A3").Value = 1
$tot.Range("B3").Value = "2022-Q3"
$tot.Range("C3").Value = 3
$tot.Range("D3").Value = 0.54

$tot.Range("A4").Value = 2
$tot.Range("B4").Value = "2022-Q2"
$tot.Range("C4").Value = 3
$tot.Range("D4").Value = 0.49

$tot.Range("A5").Value = 3
$tot.Range("B5").Value = "2022-Q1"
$tot.Range("C5").Value = 3
$tot.Range("D5").Value = 0.54

$tot.Range("A6").Value = 4
$tot.Range("B6").Value = "2021-Q3"
$tot.Range("C6").Value = 3
$tot.Range("D6").Value = 0.42

$tot.Range("A7").Value = 5
$tot.Range("B7").Value = "2021-Q2"
$tot.Range("C7").Value = 3
$tot.Range("D7").Value = 0.45
